$d = $word.ActiveDocument

# Ordered list of (old, new) text replacements, in document order.
# Using Replace:=1 (wdReplaceOne) and Wrap:=0 (wdFindStop) so each call
# only touches the single intended occurrence, searching forward from the
# start of the document content each time. This matters because one of the
# new values ("10÷4=") equals an old value earlier in the document that has
# already been replaced by the time we get to it, and replacing strictly in
# document order (top to bottom) with Wrap:=0 avoids re-matching any newly
# inserted text.
$pairs = @(
    ,@("2025-07-25 Friday", "2025-07-26 Saturday")
    ,@("27÷8=", "10÷3=")
    ,@("54÷7=", "95÷4=")
    ,@("12÷9=", "13÷7=")
    ,@("35÷2=", "28÷9=")
    ,@("98÷2=", "52÷3=")
    ,@("67÷5=", "36÷6=")
    ,@("91÷2=", "95÷9=")
    ,@("80÷9=", "97÷2=")
    ,@("10÷4=", "25÷8=")
    ,@("45÷3=", "64÷7=")
    ,@("26÷2=", "48÷2=")
    ,@("76÷5=", "95÷7=")
    ,@("43÷6=", "57÷2=")
    ,@("69÷4=", "37÷8=")
    ,@("43÷5=", "44÷7=")
    ,@("30÷9=", "30÷2=")
    ,@("45÷7=", "68÷9=")
    ,@("47÷7=", "29÷2=")
    ,@("51÷2=", "52÷4=")
    ,@("10÷8=", "10÷4=")
    ,@("40÷8=", "45÷9=")
    ,@("52÷5=", "18÷2=")
    ,@("29÷5=", "84÷7=")
    ,@("99÷2=", "72÷5=")
    ,@("81÷6=", "31÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 1) | Out-Null
}

Write-Output "Done."
